$d = $word.ActiveDocument

# Curly single quotes used around 'title' in the document text.
$lq = [char]0x2018
$rq = [char]0x2019

# ------------------------------------------------------------------
# 1) "...sorting all of the individual columns..." -> "...sorting all
#    the individual columns..." (drop "of", merging the runs/proofErr
#    markers that previously isolated "all of").
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "through each row sorting all of the individual",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "through each row sorting all the individual", 2)

# ------------------------------------------------------------------
# 2) "...We then added a country ID to each video related csv..." --
#    move "each " into the previous sentence/run, and fix the
#    "Postgresql" -> "Postgres" typo further along the same run.
# ------------------------------------------------------------------
$old2 = "each video related csv to denote the country of origin. The next major transformation came about when we ran into repeated loading failures with Postgresql. Initially when we attempted to load the data into SQL, the load would fail, citing that there was an unterminated cell in the " + $lq + "title" + $rq + " column. Upon closer inspection we realized that many of the rows contained commas in the " + $lq + "title" + $rq + " column which broke our load. To combat "
$new2 = "video related csv to denote the country of origin. The next major transformation came about when we ran into repeated loading failures with Postgres. Initially when we attempted to load the data into SQL, the load would fail, citing that there was an unterminated cell in the " + $lq + "title" + $rq + " column. Upon closer inspection we realized that many of the rows contained commas in the " + $lq + "title" + $rq + " column which broke our load. To combat "
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

$rBoundary = $d.Content
$rBoundary.Find.Execute("country ID to ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ip = $d.Range($rBoundary.End, $rBoundary.End)
$ip.InsertBefore("each ")

# ------------------------------------------------------------------
# 3) "...exported to csv to be loaded into MySQL." -> "...loaded into
#    Postgres."
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "exported to csv to be loaded into MySQL.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "exported to csv to be loaded into Postgres.", 2)

# ------------------------------------------------------------------
# 4) "...data that we collected using PostgrSQL. The reason why..."
#    -> "...using Postgres. The reason why..." (fix the other typo,
#    the MySQL-front-end sentence earlier in the doc is untouched).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "data that we collected using PostgrSQL. The reason why",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "data that we collected using Postgres. The reason why", 2)
